$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeDataBatch16")

# Update the values in column E (rows 2-4) with the new string values.
$ws.Range("E2").Value = "cold789"
$ws.Range("E3").Value = "hot8945"
$ws.Range("E4").Value = "worm357"

# Move the active cell selection from E4 to D4.
$ws.Range("D4").Select()
